$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.917.96'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '2.230.94'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.569'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.77%  '
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.57%  '
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '2.570.55'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = '2.231.29'
$ws.Range("E15").Value = '  -2.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.819'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.79%  '
$ws.Range("D18").Value = '43.785.47'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").Value = '0.0₃0958'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -10.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0795'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.22'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("E35").Value = '  -3.07%  '
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.118'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -9.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0296'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.45%  '
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("D44").Value = '1.718.77'
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '84.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.186'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.26%  '
